# Emptying the charts in the word documents
# Adds a new "At Surgery / 1 Week Post-OP" results table (with a leading
# note row) to the bottom of Sheet1, and tweaks a few header row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Adjust the height of a few existing header rows (pure formatting)
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 30.75

# ---------------------------------------------------------------------
# 2. Row 100 : a free-floating note, styled like the existing
#    "borderless note" rows (e.g. A22/A23) - Arial 12, wrap, top aligned,
#    but with no border/fill.
# ---------------------------------------------------------------------
$ws.Range("A89").Copy()
$ws.Range("A100").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A100").Borders.LineStyle = -4142 # xlNone
$ws.Range("A100").Value = "This is in the post script"
$ws.Rows.Item(100).RowHeight = 30.75

# ---------------------------------------------------------------------
# 3. Row 101 : header row for the new table
# ---------------------------------------------------------------------
$ws.Range("A94:D94").Copy()
$ws.Range("A101:D101").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A101").Value = "Variable"
$ws.Range("B101").Value = "At Surgery"
$ws.Range("C101").Value = "1 Week Post-OP"
$ws.Range("D101").Value = "Units"
$ws.Rows.Item(101).RowHeight = 32.25

# ---------------------------------------------------------------------
# 4. Row 102 : Blood Pressure data row
# ---------------------------------------------------------------------
$ws.Range("A98:D98").Copy()
$ws.Range("A102:D102").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A102").Value = "Blood Pressure"
$ws.Range("B102").Value = "148/100"
$ws.Range("C102").Value = "135/95"
$ws.Range("D102").Value = "mmHg"
$ws.Rows.Item(102).RowHeight = 15.75

# ---------------------------------------------------------------------
# 5. Row 103 : Plasma [K+] data row
# ---------------------------------------------------------------------
$ws.Range("A98:D98").Copy()
$ws.Range("A103:D103").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A103").Value = "Plasma [K+]"
$ws.Range("B103").Value = 3.6
$ws.Range("C103").Value = 4.6
$ws.Range("D103").Value = "mEq/L"
$ws.Rows.Item(103).RowHeight = 15.75

# ---------------------------------------------------------------------
# 6. Scroll / selection bookkeeping to match where the editor ended up
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("B103"), $false)
$ws.Range("B103").Select()
$ws.Application.ActiveWindow.ScrollRow = 91

$wb.Save()
